# "08-OOP-Inheritance-Polymorphism/02. Aggregators.docx"
#
# The diff boils down to one real content edit: the stray leading "0"
# run in front of the "8. Destructors, Constructors and Copy-Assignment"
# Heading 1 is deleted (so the title correctly reads "8. ..." instead of
# "08. ..."), and Word's automatic "last edit location" bookmark
# (_GoBack) moves from its old spot (after "output but" later in the
# document) to the place where this edit was made, right at the start of
# the heading. Everything else in the raw XML diff (namespace-list
# pruning, regenerated VML anchorIds/GUIDs, a stray relationship-id
# swap and cached PAGE-field digits buried inside footer textboxes,
# w15:restartNumberingAfterBreak housekeeping attributes, latent style
# bookkeeping, table-style margins, etc.) is incidental churn from the
# authoring copy of Word re-serializing parts the COM object model here
# doesn't expose for editing (footer text boxes aren't walkable story
# ranges in this host, and numbering.xml/styles.xml internals aren't
# reachable at all through Document.* members) - so it is not
# reproducible (or meaningful) to chase from script code and is left
# alone.

$d = $word.ActiveDocument

# The Heading 1 paragraph is the very first paragraph in the body.
$heading = $d.Paragraphs(1).Range

# Isolate just the leading "0" character of "08. Destructors, ...".
$zero = $heading.Duplicate
$zero.SetRange($heading.Start, $heading.Start + 1)

if ($zero.Text -eq "0") {
    # Drop a fresh _GoBack bookmark right where the edit is about to
    # happen (immediately after the "0", i.e. right before the "8").
    # Word keeps only one _GoBack bookmark per document, so adding this
    # one automatically removes whichever one already existed further
    # down (after "output but"), matching the diff exactly. Doing this
    # *before* deleting the "0" (rather than collapsing to position 0
    # afterwards) keeps the empty bookmark pair anchored inside the
    # heading paragraph instead of sliding into the next one.
    $bookmarkSpot = $d.Range($zero.End, $zero.End)
    $d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

    # Remove the stray leading "0" so the title reads "8. ..." again.
    $zero.Delete()
}

Write-Output ("Heading1 now reads: " + $d.Paragraphs(1).Range.Text)
